# Apply the "KLGeneralInformationFSIII" term 2.0.0 update:
#  - bump Version / Date / Contact on the Metadata sheet
#  - add a second "Include from FSIII 2" sheet (a copy of "Include from FSIII")
#    whose "Value" column for the descendent-of row is recoded to a concept id
#    GUID on the first (original) sheet, and contains the literal code "F" on
#    the newly added sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# ---------------------------------------------------------------------------
# 2. Recode the descendent-of concept value on the existing
#    "Include from FSIII" sheet from "F" to the new concept GUID.
# ---------------------------------------------------------------------------
$include1 = $wb.Worksheets.Item("Include from FSIII")
$include1.Range("C2").Value = "370e6178-9a5d-45f9-a2c9-f658186059c9"

# ---------------------------------------------------------------------------
# 3. Add the new "Include from FSIII 2" sheet (copy of the first, placed
#    right after it) and give it the original "F" concept value.
# ---------------------------------------------------------------------------
$include1.Copy($null, $include1)
$include2 = $wb.Worksheets.Item($include1.Index + 1)
$include2.Name = "Include from FSIII 2"
$include2.Range("C2").Value = "F"

# Keep the originally active sheet ("Metadata") selected, as in the source file.
$meta.Activate()
